$wb = $excel.ActiveWorkbook

# Sheet ALC, row 12 (item id 5515)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 499.4
$ws.Range("I12").Value = 499.4
$ws.Range("K12").Value = 499.4
$ws.Range("M12").Value = -329.4

# Sheet ALC, row 40 (item id 5505)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2789.0908
$ws.Range("J40").Value = 4180.1665
$ws.Range("L40").Value = 4180.1665
$ws.Range("N40").Value = -4530.1665

# Sheet ALC, row 43 (item id 5472)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3055.524
$ws.Range("I43").Value = 3694.2222
$ws.Range("J43").Value = 2576.5
$ws.Range("K43").Value = 3694.2222
$ws.Range("L43").Value = 2576.5
$ws.Range("M43").Value = -3625.2222
$ws.Range("N43").Value = -2714.5

# Sheet ALC, row 53 (item id 5479)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 537.61536
$ws.Range("I53").Value = 531.5
$ws.Range("K53").Value = 531.5
$ws.Range("M53").Value = 105.5

# Sheet ALC, row 74 (item id 5507)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3772.875
$ws.Range("I74").Value = 2998
$ws.Range("J74").Value = 6097.5
$ws.Range("K74").Value = 2998
$ws.Range("L74").Value = 6097.5
$ws.Range("M74").Value = -2062
$ws.Range("N74").Value = -7969.5

# Sheet ALC, row 77 (item id 5507)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3772.875
$ws.Range("I77").Value = 2998
$ws.Range("J77").Value = 6097.5
$ws.Range("K77").Value = 14990
$ws.Range("L77").Value = 30487.5
$ws.Range("M77").Value = -10310
$ws.Range("N77").Value = -39847.5

# Sheet ALC, row 96 (item id 19894)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1493.125
$ws.Range("I96").Value = 1657.5
$ws.Range("K96").Value = 4972.5
$ws.Range("M96").Value = -3599.5

# Sheet ALC, row 135 (item id 44047)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 468.86667
$ws.Range("I135").Value = 464.15384
$ws.Range("K135").Value = 4177.38456
$ws.Range("M135").Value = -1642.38456

# Sheet ALC, row 137 (item id 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2950.2273
$ws.Range("I137").Value = 2976.9167
$ws.Range("J137").Value = 2918.2
$ws.Range("K137").Value = 8930.750100000001
$ws.Range("L137").Value = 8754.599999999999
$ws.Range("M137").Value = -6380.750100000001
$ws.Range("N137").Value = -13854.6

# Sheet ALC, row 138 (item id 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2404.68
$ws.Range("J138").Value = 3074.923
$ws.Range("L138").Value = 9224.769
$ws.Range("N138").Value = -19504.769

# Sheet ALC, row 141 (item id 44161)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 9285.817999999999
$ws.Range("I141").Value = 8404.9
$ws.Range("K141").Value = 25214.7
$ws.Range("M141").Value = -20034.7

# Sheet ARM, row 6 (item id 2226)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 501497
$ws.Range("I6").Value = 1001499
$ws.Range("J6").Value = 1495
$ws.Range("K6").Value = 1001499
$ws.Range("L6").Value = 1495
$ws.Range("M6").Value = -1001326
$ws.Range("N6").Value = -1841

# Sheet ARM, row 74 (item id 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2733.25
$ws.Range("I74").Value = 2877.4119
$ws.Range("K74").Value = 2877.4119
$ws.Range("M74").Value = -2003.4119

# Sheet ARM, row 77 (item id 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2733.25
$ws.Range("I77").Value = 2877.4119
$ws.Range("K77").Value = 14387.0595
$ws.Range("M77").Value = -10019.0595

# Sheet ARM, row 80 (item id 10667)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 26566.143
$ws.Range("J80").Value = 29992.8
$ws.Range("L80").Value = 29992.8
$ws.Range("N80").Value = -31988.8

# Sheet ARM, row 83 (item id 10667)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 26566.143
$ws.Range("J83").Value = 29992.8
$ws.Range("L83").Value = 89978.39999999999
$ws.Range("N83").Value = -99962.39999999999

# Sheet ARM, row 105 (item id 18699)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 22174.5
$ws.Range("J105").Value = 22174.5
$ws.Range("L105").Value = 22174.5
$ws.Range("N105").Value = -29162.5

# Sheet BSM, row 103 (item id 18514)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""

# Sheet BSM, row 105 (item id 19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4246.077
$ws.Range("I105").Value = 2424.75
$ws.Range("J105").Value = 5055.5557
$ws.Range("K105").Value = 2424.75
$ws.Range("L105").Value = 5055.5557
$ws.Range("M105").Value = -677.75
$ws.Range("N105").Value = -8549.555700000001

# Sheet BSM, row 107 (item id 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1797.8
$ws.Range("J107").Value = 1000
$ws.Range("L107").Value = 1000
$ws.Range("N107").Value = -4840

# Sheet CRP, row 2 (item id 1820)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1107.5
$ws.Range("I2").Value = 765.7143
$ws.Range("J2").Value = 3500
$ws.Range("K2").Value = 765.7143
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = -652.7143
$ws.Range("N2").Value = -3726

# Sheet CRP, row 9 (item id 15611)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 169000
$ws.Range("J9").Value = 169000
$ws.Range("L9").Value = 169000
$ws.Range("N9").Value = -169336

# Sheet CRP, row 10 (item id 1997)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 399.6
$ws.Range("J10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("N10").Value = -1278

# Sheet CRP, row 31 (item id 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2270.4285
$ws.Range("I31").Value = 4999
$ws.Range("J31").Value = 1179
$ws.Range("K31").Value = 4999
$ws.Range("L31").Value = 1179
$ws.Range("M31").Value = -4704
$ws.Range("N31").Value = -1769

# Sheet CRP, row 34 (item id 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2270.4285
$ws.Range("I34").Value = 4999
$ws.Range("J34").Value = 1179
$ws.Range("K34").Value = 4999
$ws.Range("L34").Value = 1179
$ws.Range("M34").Value = -4797
$ws.Range("N34").Value = -1583

# Sheet CUL, row 4 (item id 4650)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 722022.6
$ws.Range("I4").Value = 722022.6
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2166067.8
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -2165955.8
$ws.Range("N4").Value = ""

# Sheet CUL, row 9 (item id 4681)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 6622.375
$ws.Range("I9").Value = 993.3333
$ws.Range("K9").Value = 2979.9999
$ws.Range("M9").Value = -2755.9999

# Sheet CUL, row 12 (item id 4854)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 123.57143
$ws.Range("J12").Value = 4.2
$ws.Range("L12").Value = 12.6
$ws.Range("N12").Value = -358.6

# Sheet CUL, row 23 (item id 4858)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 159.28572
$ws.Range("I23").Value = 144.83333
$ws.Range("J23").Value = 246
$ws.Range("K23").Value = 434.49999
$ws.Range("L23").Value = 738
$ws.Range("M23").Value = -199.49999
$ws.Range("N23").Value = -1208

# Sheet CUL, row 34 (item id 4749)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3900.0908
$ws.Range("J34").Value = 4225.2
$ws.Range("L34").Value = 12675.6
$ws.Range("N34").Value = -12843.6

# Sheet CUL, row 39 (item id 4712)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 9000
$ws.Range("J39").Value = 9000
$ws.Range("L39").Value = 27000
$ws.Range("N39").Value = -27588

# Sheet CUL, row 106 (item id 19819)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 8993.333000000001
$ws.Range("I106").Value = 9980
$ws.Range("J106").Value = 8500
$ws.Range("K106").Value = 29940
$ws.Range("L106").Value = 25500
$ws.Range("M106").Value = -28994
$ws.Range("N106").Value = -27392

# Sheet CUL, row 114 (item id 27865)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 499
$ws.Range("I114").Value = 499.5
$ws.Range("J114").Value = 498
$ws.Range("K114").Value = 1498.5
$ws.Range("L114").Value = 1494
$ws.Range("M114").Value = 1755.5
$ws.Range("N114").Value = -8002

# Sheet CUL, row 115 (item id 27861)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 300
$ws.Range("I115").Value = 300
$ws.Range("K115").Value = 900
$ws.Range("M115").Value = 275

# Sheet CUL, row 139 (item id 44102)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1332.75
$ws.Range("I139").Value = 1332.75
$ws.Range("K139").Value = 3998.25
$ws.Range("M139").Value = 1141.75

# Sheet GSM, row 97 (item id 19940)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 388.6
$ws.Range("J97").Value = 573.2
$ws.Range("L97").Value = 573.2
$ws.Range("N97").Value = -1565.2

# Sheet GSM, row 122 (item id 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2480.182
$ws.Range("I122").Value = 2518.2
$ws.Range("K122").Value = 7554.599999999999
$ws.Range("M122").Value = -5104.599999999999

# Sheet GSM, row 126 (item id 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2052
$ws.Range("I126").Value = 2123.111
$ws.Range("J126").Value = 1732
$ws.Range("K126").Value = 6369.333
$ws.Range("L126").Value = 5196
$ws.Range("M126").Value = -3899.333
$ws.Range("N126").Value = -10136

# Sheet GSM, row 132 (item id 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6998.5
$ws.Range("I132").Value = 6998.5
$ws.Range("K132").Value = 20995.5
$ws.Range("M132").Value = -18465.5

# Sheet LTW, row 4 (item id 3788)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = ""

# Sheet LTW, row 16 (item id 5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 721.1667
$ws.Range("I16").Value = 665.4
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 665.4
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -495.4
$ws.Range("N16").Value = -1340

# Sheet LTW, row 28 (item id 3788)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").Value = ""

# Sheet LTW, row 37 (item id 3788)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").Value = ""

# Sheet LTW, row 55 (item id 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1138.4667
$ws.Range("I55").Value = 702.25
$ws.Range("J55").Value = 1637
$ws.Range("K55").Value = 702.25
$ws.Range("L55").Value = 1637
$ws.Range("M55").Value = -529.25
$ws.Range("N55").Value = -1983

# Sheet LTW, row 69 (item id 10671)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 28999.5
$ws.Range("J69").Value = 28999.5
$ws.Range("L69").Value = 28999.5
$ws.Range("N69").Value = -30621.5

# Sheet LTW, row 72 (item id 10671)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H72").Value = 28999.5
$ws.Range("J72").Value = 28999.5
$ws.Range("L72").Value = 86998.5
$ws.Range("N72").Value = -95110.5

# Sheet LTW, row 115 (item id 26015)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 30151
$ws.Range("J115").Value = 30151
$ws.Range("L115").Value = 30151
$ws.Range("N115").Value = -32501

# Sheet WVR, row 29 (item id 3568)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 8950.666999999999
$ws.Range("I29").Value = 1582.5454
$ws.Range("J29").Value = 90000
$ws.Range("K29").Value = 1582.5454
$ws.Range("L29").Value = 90000
$ws.Range("M29").Value = -1292.5454
$ws.Range("N29").Value = -90580

# Sheet WVR, row 126 (item id 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1964.75
$ws.Range("I126").Value = 1987.5714
$ws.Range("J126").Value = 1805
$ws.Range("K126").Value = 5962.7142
$ws.Range("L126").Value = 5415
$ws.Range("M126").Value = -3492.7142
$ws.Range("N126").Value = -10355

# Sheet WVR, row 136 (item id 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3844.4194
$ws.Range("I136").Value = 4150.1924
$ws.Range("J136").Value = 2254.4
$ws.Range("K136").Value = 12450.5772
$ws.Range("L136").Value = 6763.200000000001
$ws.Range("M136").Value = -9900.5772
$ws.Range("N136").Value = -11863.2
